$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect the unified "DataNode" concept
# (was "Property1") per the commit: unify DataNode / DataTable / Entity.
$ws.Name = "DataNode"

# Preserve the editor's last selection on the sheet.
$ws.Range("C41").Select()
